$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 52, shifting existing rows 52-112 down to 53-113
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new weekly data point
$ws.Range("A52").Value = 6
$ws.Range("B52").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C52").Value = "Metropolitana"
$ws.Range("D52").Value = 44494
$ws.Range("E52").Value = 13
$ws.Range("F52").Value = 100112001
$ws.Range("G52").Value = "Berenjena"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 580
$ws.Range("K52").Value = 6000
$ws.Range("L52").Value = 7000
$ws.Range("M52").Value = 6448
$ws.Range("N52").Value = "`$/caja 50 unidades"
$ws.Range("O52").Value = "Región de Arica y Parinacota"
$ws.Range("P52").Value = 129
$ws.Range("Q52").Value = 50
$ws.Range("R52").Value = "Hortaliza"
